$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 3 (shifts existing rows 3..73 down to 4..74),
#    creating room for the new "B0_P3_fim_ano_lectivo" field right after
#    "B0_P2_inicio_ano_lectivo".
$ws.Rows.Item(3).Insert()

# 2) Update the constraint/message for B0_P2_inicio_ano_lectivo (row 2):
#    the allowed year range changed from "2024 or 2025" to a 2000-2024 regex.
$ws.Range("C2").Value = 'regex(., ''^(200[0-9]|201[0-9]|202[0-4])$'')'
$ws.Range("D2").Value = 'O ano deve estar entre 2000 e 2024, e não pode começar com 1'

# 3) Fill in the newly inserted row 3 with the new field
#    "B0_P3_fim_ano_lectivo" and its constraint/message.
$ws.Range("A3").Value = 'B0_P3_fim_ano_lectivo'
$ws.Range("C3").Value = '. = ${QEPE_DGE_SQE_B0_P2_inicio_ano_lectivo} or . = ${QEPE_DGE_SQE_B0_P2_inicio_ano_lectivo} + 1'
$ws.Range("D3").Value = 'O fim do ano letivo deve ser igual ao início ou ao ano seguinte'

# 4) Several text fields now require the value to start with a letter.
#    Their validation message stays the same; only the constraint (column C)
#    changes. After the row insertion above these fields live in rows 4-10.
$startsWithLetter = 'regex(., ''^[A-Za-z].*'')'
$ws.Range("C4").Value = $startsWithLetter   # nome_escola
$ws.Range("C5").Value = $startsWithLetter   # endereco_escola
$ws.Range("C6").Value = $startsWithLetter   # ponto_referencia
$ws.Range("C7").Value = $startsWithLetter   # localidade
$ws.Range("C8").Value = $startsWithLetter   # decreto_criacao
$ws.Range("C9").Value = $startsWithLetter   # licenca
$ws.Range("C10").Value = $startsWithLetter  # nome

# 5) area_formacao (now row 12) is restricted to letters/accented
#    letters/spaces only, with a new validation message.
$ws.Range("C12").Value = 'regex(., ''^[A-Za-zÀ-ÿ\s]+$'')'
$ws.Range("D12").Value = 'Só são permitidas letras e espaços. Números e símbolos não são aceitos'

# 6) Update the hidden filter-database defined name to match the new
#    used range (one extra row).
foreach ($dn in $wb.Names) {
    if ($dn.Name -like "*_FilterDatabase*") {
        $dn.RefersTo = "=Plan1!`$A`$1:`$D`$74"
    }
}

# 7) Update the sheet view state (scroll position / zoom / selection) to
#    match the saved view from the edited workbook.
$ws.Range("C25").Select()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 9
